$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of data (rows 6-8) beneath the existing table.
$data = @(
    @(1, 105, 8.5, $true),
    @(2, 105, 9.5, $true),
    @(3, 105, 9,   $true)
)

$r = 6
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Match the formatting already applied to the existing data rows (A2:D5):
# vertical-centered, wrapped text (style index 2 in the workbook's cellXfs).
$newRange = $ws.Range("A6:D8")
$newRange.VerticalAlignment = -4108
$newRange.WrapText = $true

# Reflect the final cursor/selection position left by data entry.
$ws.Range("C8").Select()
